$wb = $excel.ActiveWorkbook

# Rename the "Include" sheet.
$includeSheet = $wb.Worksheets.Item("Include from Ferlab.bio CodeS")
$includeSheet.Name = "Include #0"

$meta = $wb.Worksheets.Item("Metadata")

# New publish date.
$meta.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# New contact value.
$meta.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# Make room for the new "Jurisdiction" row (shifts rows 11-14 down to 12-15).
$meta.Range("A11:B11").Insert()

# Copy the style of the row that landed right below (now "Description")
# onto the blank new row so it matches the rest of the table's
# borders/alignment.
$meta.Range("A12:B12").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)

$meta.Range("A11").Value = "Jurisdiction"

# B11 should be an explicit empty string (matching the "no display value"
# convention already used elsewhere in this workbook), not merely a blank
# cell - paste an existing empty-string cell's value into it.
$includeSheet.Range("A3").Copy()
$meta.Range("B11").PasteSpecial(-4163)
